$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price / column D value or $null if unchanged, new Volume(1h) / column E value)
$changes = @{
    2  = @("60.768.50", "  +0.05%  ")
    3  = @("3.364.74",  "  -0.56%  ")
    4  = @("0.999",     "  -0.02%  ")
    5  = @("569.49",    "  +0.17%  ")
    6  = @("137.69",    "  -2.04%  ")
    7  = @($null,       "  -0.06%  ")
    8  = @($null,       "  -0.67%  ")
    9  = @("7.69",      "  +2.84%  ")
    10 = @($null,       "  -1.85%  ")
    11 = @($null,       "  -4.12%  ")
    12 = @("3.938.13",  "  -0.58%  ")
    13 = @($null,       "  +0.55%  ")
    14 = @("27.88",     "  -2.00%  ")
    15 = @("3.366.71",  "  -0.55%  ")
    16 = @($null,       "  -1.80%  ")
    17 = @("60.875.03", "  +0.08%  ")
    18 = @($null,       "  -1.82%  ")
    19 = @("13.49",     "  -3.47%  ")
    20 = @("8.89",      "  -0.87%  ")
    21 = @("382.58",    "  -0.25%  ")
    22 = @("75.18",     "  +2.14%  ")
    23 = @($null,       "  -2.07%  ")
    24 = @($null,       "  -0.20%  ")
    25 = @($null,       "  -5.85%  ")
    26 = @($null,       "  +6.76%  ")
    27 = @($null,       "  +0.07%  ")
    28 = @("7.13",      "  -4.03%  ")
    29 = @("7.84",      "  -1.67%  ")
    30 = @($null,       "  -1.79%  ")
    31 = @($null,       "  -0.05%  ")
    32 = @($null,       "  -7.21%  ")
    33 = @("22.96",     "  -2.85%  ")
    34 = @("167.31",    "  +0.37%  ")
    35 = @($null,       "  -1.83%  ")
    36 = @("4.91",      "  -1.70%  ")
    37 = @("3.398.34",  "  -0.43%  ")
    38 = @("1.44",      "  -3.14%  ")
    39 = @("0.0754",    "  -2.66%  ")
    40 = @("25.34",     "  -9.00%  ")
    41 = @("0.769",     "  -1.38%  ")
    42 = @("4.33",      "  -1.88%  ")
    43 = @("1.62",      "  -3.13%  ")
    44 = @($null,       "  -1.79%  ")
    45 = @("2.444.39",  "  -2.73%  ")
    46 = @("1.00",      "  -0.02%  ")
    47 = @("6.62",      "  -3.10%  ")
    48 = @("22.07",     "  -6.13%  ")
    49 = @("0.0256",    "  -5.02%  ")
    50 = @("1.93",      "  -6.77%  ")
    51 = @("0.202",     "  -2.79%  ")
}

foreach ($row in $changes.Keys) {
    $pair = $changes[$row]
    $priceValue = $pair[0]
    $volumeValue = $pair[1]

    if ($priceValue -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $priceValue
    }
    $ecell = $ws.Cells.Item($row, 5)
    $ecell.NumberFormat = "@"
    $ecell.Value = $volumeValue
}
